# Generate Report for Handback
#
# The localization round-trip for file "290d4694-5fbe-482f-b96d-8629f8316313"
# (both the zh-cn and de-de targets) has come back from the vendor "in sync
# with en-US". Update the Overview sheet and the per-locale detail sheets to
# reflect the handback: status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Target File / Latest Handback
# File / Latest Handback DateTime columns get populated, and a hyperlink is
# added for the newly-populated target-file cell.
#
# The other tracked file, "4f5603d3-558f-44f4-8725-eee2aa8e07d3", is
# untouched by this handback and keeps its "Ready for handoff" status.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$targetFileName = "290d4694-5fbe-482f-b96d-8629f8316313.md"
$targetFileUrl  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/0ca81370325b662d381c4186fd6f48056150f798/e2e/290d4694-5fbe-482f-b96d-8629f8316313.md"

$zhHandbackFile = "290d4694-5fbe-482f-b96d-8629f8316313.ffc62124b4d2bc2f5e6f0423e7c87925ef22769b.zh-cn.xlf"
$deHandbackFile = "290d4694-5fbe-482f-b96d-8629f8316313.ffc62124b4d2bc2f5e6f0423e7c87925ef22769b.de-de.xlf"

$zhHandbackDateTime = "2016-08-14 02:56:15"
$deHandbackDateTime = "2016-08-14 02:56:25"

# ---------------------------------------------------------------------------
# Overview sheet: update the per-locale status columns (E = zh-cn, F = de-de)
# for the row belonging to 290d4694-5fbe-482f-b96d-8629f8316313.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn detail sheet: update Status, Latest Target File (+hyperlink),
# Latest Handback File and Latest Handback DateTime for row 2.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetFileUrl, [Type]::Missing, [Type]::Missing, $targetFileName)
$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("K2").Value = $zhHandbackDateTime

# ---------------------------------------------------------------------------
# de-de detail sheet: same updates as zh-cn, with the de-de handback file
# and its own handback datetime.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetFileUrl, [Type]::Missing, [Type]::Missing, $targetFileName)
$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("K2").Value = $deHandbackDateTime
